$wb = $excel.ActiveWorkbook

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md"
$wsOverview.Range("A3").Value = "ffffff5a4218d1-3aac-450e-954f-2c04d7354d70.md"
$wsOverview.Range("A4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md"
$wsZhCn.Range("C2").Value = "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-03-01 09:38:21"
$wsZhCn.Range("E2").Value = "b6091237-6809-4684-867b-5538749eeb17.md"
$wsZhCn.Range("F2").Value = "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf"
$wsZhCn.Range("G2").Value = "2016-03-01 09:39:07"

$wsZhCn.Range("A3").Value = "ffffff5a4218d1-3aac-450e-954f-2c04d7354d70.md"
$wsZhCn.Range("C3").Value = "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-01 09:38:21"
$wsZhCn.Range("E3").Value = "b6091237-6809-4684-867b-5538749eeb17.md"
$wsZhCn.Range("F3").Value = "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.zh-cn.xlf"
$wsZhCn.Range("G3").Value = "2016-03-01 09:39:07"

$wsZhCn.Range("A4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.md"
$wsZhCn.Range("B4").Value = "Ready for handoff"
$wsZhCn.Range("C4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.372ffa13721dd814a514de1f72dd7f1ee4531b68.zh-cn.xlf"
$wsZhCn.Range("D4").Value = "2016-03-01 09:44:23"
$wsZhCn.Range("E4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.md"
$wsZhCn.Range("F4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.372ffa13721dd814a514de1f72dd7f1ee4531b68.zh-cn.xlf"
$wsZhCn.Range("G4").Value = "2016-03-01 09:43:24"

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "ffff0ac76892-4552-4a0c-bfb9-ea3d72729ca7.md"
$wsDeDe.Range("C2").Value = "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-03-01 09:38:31"
$wsDeDe.Range("E2").Value = "b6091237-6809-4684-867b-5538749eeb17.md"
$wsDeDe.Range("F2").Value = "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf"
$wsDeDe.Range("G2").Value = "2016-03-01 09:39:25"

$wsDeDe.Range("A3").Value = "ffffff5a4218d1-3aac-450e-954f-2c04d7354d70.md"
$wsDeDe.Range("C3").Value = "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-01 09:38:31"
$wsDeDe.Range("E3").Value = "b6091237-6809-4684-867b-5538749eeb17.md"
$wsDeDe.Range("F3").Value = "b6091237-6809-4684-867b-5538749eeb17.850ce1640dff362fa460518f8d850d8796ed11b7.de-de.xlf"
$wsDeDe.Range("G3").Value = "2016-03-01 09:39:25"

$wsDeDe.Range("A4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.md"
$wsDeDe.Range("B4").Value = "Ready for handoff"
$wsDeDe.Range("C4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.372ffa13721dd814a514de1f72dd7f1ee4531b68.de-de.xlf"
$wsDeDe.Range("D4").Value = "2016-03-01 09:44:34"
$wsDeDe.Range("E4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.md"
$wsDeDe.Range("F4").Value = "e4e33f3a-4f22-481e-bf49-8a7ddc7d3734.372ffa13721dd814a514de1f72dd7f1ee4531b68.de-de.xlf"
$wsDeDe.Range("G4").Value = "2016-03-01 09:43:42"
